# PRM16: add a embedding to prm15_2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: clear the "Running" status label in A24 (status cell becomes blank,
#     keeps its existing formatting) ---
$ws.Range("A24").ClearContents()

# --- Row 27: populate the new PRM15_2 result row ---
# A27 mirrors the "Running" status used elsewhere (blue status style), so copy
# formatting+value from A22 (same value/style) rather than retyping it.
$ws.Range("A22").Copy($ws.Range("A27"))

# K27 keeps the same (empty-but-highlighted) style as K26 but has no value.
$ws.Range("K26").Copy($ws.Range("K27"))
$ws.Range("K27").ClearContents()

# Remaining cells: plain values
$ws.Range("B27").Value = "100/30"
$ws.Range("C27").Value = "no"
$ws.Range("D27").Value = 64
$ws.Range("E27").Value = "dotproduct"
$ws.Range("F27").Value = "N"
$ws.Range("G27").Value = "Y"
$ws.Range("H27").Value = 16
$ws.Range("I27").Value = "temp1"
$ws.Range("J27").Value = 64
# Note: M27 is written before L27 so the two brand-new shared strings land in
# the same order as the authored workbook (description first, then the name).
$ws.Range("M27").Value = "based on SGE, (mean+max)/2, dotproduct/sqrt(c/g)"
$ws.Range("L27").Value = "prm15_2_resnet50"

# --- Selection / scroll position ---
[void]$ws.Range("F27").Select()
